# This script reproduces the OOXML diff describing the edit to
# "SE-N17-047.docx":
#  - remove the (stray) <w:tblPrEx> table-row-property-exception blocks
#    that precede every <w:tr> in the document's single table
#  - drop the "_GoBack" bookmark that sat after "SE-N17-047" in the
#    "Test ID" cell
#  - change "Proof of testing:" to "Proof of testing" and remove the two
#    blank paragraphs that immediately followed it
#  - append explanatory text to the "Comments: " paragraph (the group
#    did not attempt the server-to-server requirements, so the test is
#    recorded as a fail), and move the "_GoBack" bookmark to sit after
#    that new text
#  - remove the stray <w:lastRenderedPageBreak/> before
#    "Description of changes made:"
#
# <w:tblPrEx> is not reachable through the normal Word table/row/cell
# object model (Word never exposes it as a settable property on Row),
# so the only reliable way to drop it is to rebuild the table's XML
# (folding in the other textual edits above) and splice it back into
# the document with Range.InsertXML after removing the old table.

$d = $word.ActiveDocument

$tableXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblW w:w="9638" w:type="dxa"/><w:tblInd w:w="45" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tblBorders><w:tblCellMar><w:left w:w="10" w:type="dxa"/><w:right w:w="10" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3358"/><w:gridCol w:w="6280"/></w:tblGrid><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidTr="005A1147"><w:tc><w:tcPr><w:tcW w:w="3358" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>TEST LOG FORM</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6280" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="0015082F"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Test Log No: </w:t></w:r><w:r w:rsidR="002D5DD0"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>47</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidTr="005A1147"><w:tc><w:tcPr><w:tcW w:w="3358" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Test ID: </w:t></w:r><w:r w:rsidR="002D5DD0"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>SE-N17-047</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6280" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Test Date: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>30/01/2013</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidTr="005A1147"><w:tc><w:tcPr><w:tcW w:w="3358" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="00AF317F"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Tester: </w:t></w:r><w:r w:rsidR="00AF317F"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Matt</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6280" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Group: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>17</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147"><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Version When Tested</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>:</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147"><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="0015082F" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Test Passed? (Y/N): </w:t></w:r></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Passed Testing </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> ],2</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>[ ],3</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>[ ],4</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>[ ] time.</w:t></w:r></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Other…………</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00F764A8" w:rsidRPr="005A1147"><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Description of test:</w:t></w:r></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00F764A8" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:before="113" w:after="113" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009A40AC" w:rsidRPr="005A1147"><w:trPr><w:trHeight w:val="2875"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="009A40AC" w:rsidRDefault="009A40AC" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Proof of testing</w:t></w:r></w:p><w:p w:rsidR="009A40AC" w:rsidRDefault="009A40AC" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="009A40AC" w:rsidRDefault="009A40AC" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="009A40AC" w:rsidRDefault="009A40AC" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="009A40AC" w:rsidRPr="005A1147" w:rsidRDefault="009A40AC" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="005A1147" w:rsidRPr="005A1147"><w:trPr><w:trHeight w:val="2875"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Comments: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Due to the fact we did not attempt to complete the server to server requirements, there was nothing to test, so it must be regarded as a fail. This is due to time constraints and the lack of another group also being unable to achieve this task.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr><w:tr w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidTr="005A1147"><w:trPr><w:trHeight w:val="518"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="FF0000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="FF0000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>The following section only applies if the test failed and changes needed to be made. If change need to be made a change control form needs to be completed documenting this. The details of that form also need to be highlighted below.</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidTr="005A1147"><w:trPr><w:trHeight w:val="518"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9638" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:left w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="2" w:space="0" w:color="000000"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tcMar><w:top w:w="55" w:type="dxa"/><w:left w:w="55" w:type="dxa"/><w:bottom w:w="55" w:type="dxa"/><w:right w:w="55" w:type="dxa"/></w:tcMar></w:tcPr><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Change Control Number:</w:t></w:r></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005A1147"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Description of changes made:</w:t></w:r></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A1147" w:rsidRPr="005A1147" w:rsidRDefault="005A1147" w:rsidP="005A1147"><w:pPr><w:pStyle w:val="Default"/><w:snapToGrid w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:textAlignment w:val="top"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$table = $d.Tables.Item(1)
$table.Delete() | Out-Null

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($tableXml) | Out-Null

Write-Output "Table rebuilt; tables now: $($d.Tables.Count)"
